# "Generate Report for Handoff"
#
# The status moves from "In Translation" to "Ready for handoff" and the
# two "Latest ... Datetime" timestamps that are recorded alongside it are
# refreshed to the moment the handoff report was produced. The "Status"
# column is also widened a bit on every sheet so the new (longer) status
# text isn't truncated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
# Overview sheet: E2 = zh-cn status, F2 = de-de status
$wsOverview.Range("E2").Value2 = "Ready for handoff"
$wsOverview.Range("F2").Value2 = "Ready for handoff"
# Per-locale detail sheets: C2 = Status
$wsZhCn.Range("C2").Value2 = "Ready for handoff"
$wsDeDe.Range("C2").Value2 = "Ready for handoff"

# --- Timestamps refreshed by the handoff report -----------------------------
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2 "Latest Handoff
# Datetime" both carried the same timestamp and both move to the new one.
$wsOverview.Range("G2").Value2 = "2016-08-24 11:01:01"
$wsDeDe.Range("H2").Value2 = "2016-08-24 11:01:01"

# zh-cn!H2 "Latest Handoff Datetime" moves to its own refreshed timestamp.
$wsZhCn.Range("H2").Value2 = "2016-08-24 11:00:55"

# --- Widen the "Status" column so "Ready for handoff" fits ------------------
# Overview: column E (zh-cn status) and column F (de-de status)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
# zh-cn / de-de detail sheets: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
